$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 are cleared (values removed); C2 and E2 get new values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -8.8770339292687108
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -6.2868300483710104

# Row 3: update values in B3:E3
$ws.Range("B3").Value = -12.284955612774002
$ws.Range("C3").Value = 2.331077913522916
$ws.Range("D3").Value = -10.224154572232669
$ws.Range("E3").Value = 21.125455078223855

# Update the active selection to match new range B1:E3
$ws.Range("B1:E3").Select()
